$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 currently holds the record dated 2021-12-27 (serial 44557) with
# prices 13000/14000/13500/750. That whole record needs to be preserved as
# a new row 5, while row 4 is updated to a newer weekly record
# (2023-01-31 / serial 44957) with new prices.

# 1) Write the preserved (old) record into row 5, copying straight from row 4.
for ($col = 1; $col -le 18; $col++) {
    $src = $ws.Cells.Item(4, $col)
    $dst = $ws.Cells.Item(5, $col)
    $dst.Value2 = $src.Value2
}
# Column D (4) carries the date number format in this sheet - match it.
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat

# 2) Update row 4 with the new weekly values.
$ws.Range("D4").Value = 44957
$ws.Range("K4").Value = 21000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 21500
$ws.Range("P4").Value = 1194

$wb.Save()
